$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C values (rows 3-15) from 0.105 to 0.0798
for ($r = 3; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = 0.079799999999999996
}

# Update the active cell selection to J10
$ws.Range("J10").Select()
